$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.24002744690420741,
    -0.0099999987598664575,
    -0.0089999987555042793,
    0.061988487830195282,
    -0.005999998785518379,
    -0.0059999987632366469,
    -0.01999999856293222,
    -0.019999998555569221,
    -0.0059999987471863747,
    -0.0059999987431496038,
    -0.0044999987641816119,
    -0.0059999987426446744,
    0.054437725003959159,
    -0.01199999865125978,
    -0.019116548861301652,
    -0.0059999987260574983,
    -0.0059999987212204786,
    -0.0089999986779991659,
    -0.0089999987723698993,
    -0.0089999987362290312,
    -0.0089999987286750738,
    -0.062351734190375208,
    -0.0089999987442466178,
    -0.041999998258022586,
    -0.041999998248549275,
    -0.0059999987626753182,
    -0.0059999987618954975,
    -0.0059999987584768988,
    -0.01199999867276702,
    -0.019999998560271237,
    0.062494406107528988,
    -0.030264408750954708,
    -0.0059999987619150374
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
